# Update cryptos list values (prices & 1h volume change %) per diff.
# Price cells (column D) are written with a leading apostrophe to force
# text storage (matches original inlineStr text cells), then ClearFormats()
# removes the transient "quote prefix" style so no stray cell style is left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.052.95"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "'3.504.67"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'584.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").Value = "'132.23"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("D7").Value = "'3.501.37"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").Value = "'7.20"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").Value = "'4.100.90"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "'27.65"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "'3.485.62"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").Value = "'64.104.66"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("E19").Value = "  +5.35%  "
$ws.Range("D20").Value = "'14.41"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "'5.69"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "'387.23"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").Value = "'0.579"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "'3.640.81"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "'73.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "'0.0000114"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D30").Value = "'0.995"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").Value = "'8.32"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").Value = "'3.507.44"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "'23.80"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'5.34"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("D38").Value = "'1.59"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("D39").Value = "'6.95"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("D40").Value = "'163.86"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.52%  "
$ws.Range("D41").Value = "'0.0804"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.95%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").Value = "'1.23"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("D48").Value = "'1.64"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").Value = "'6.89"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "'2.440.07"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.23%  "
$ws.Range("D51").Value = "'0.895"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.04%  "
# Row 28/29 swap: RenderToken/Fetch.AI -> Fetch.AI/RenderToken with updated values
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "'1.56"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.80%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.49"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.66%  "
# Row 42/43 swap: EnergySwap/Mantle -> Mantle/EnergySwap with updated values
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'0.816"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'26.31"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.42%  "
